$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results_of_testing_different_models")

$ws.Range("C2").Value = 3.307767391204834
$ws.Range("C3").Value = 3.389590978622437
$ws.Range("C4").Value = 5.408395290374756
$ws.Range("C5").Value = 6.40941309928894
$ws.Range("C6").Value = 3.98141622543335
$ws.Range("C7").Value = 4.376898050308228
